# Update latest output (run 131)
$wb = $excel.ActiveWorkbook

$wsSchedule = $wb.Worksheets.Item("Schedule")
$wsDetailed = $wb.Worksheets.Item("Detailed")

# --- Schedule sheet: rows 4 and 5 (Cost / Unit Cost columns E/F) ---
$wsSchedule.Range("E4").Value = 578.092008
$wsSchedule.Range("F4").Value = 30.58687873015873
$wsSchedule.Range("E5").Value = -205.96203225
$wsSchedule.Range("F5").Value = -6.05414556878307

# --- Detailed sheet: Price column (B) updates, plus two Type (C) label flips ---
$wsDetailed.Range("B39").Value = 33.01475
$wsDetailed.Range("B40").Value = 56.98
$wsDetailed.Range("B41").Value = 75.49753
$wsDetailed.Range("C41").Value = "historical"
$wsDetailed.Range("B42").Value = 84.41579
$wsDetailed.Range("C42").Value = "historical"
$wsDetailed.Range("B43").Value = 79.13723
$wsDetailed.Range("B44").Value = 75.11698
$wsDetailed.Range("B45").Value = 63.01881
$wsDetailed.Range("B46").Value = 60.56064
$wsDetailed.Range("B47").Value = 65
$wsDetailed.Range("B48").Value = 63.62629
$wsDetailed.Range("B49").Value = 61.36777

$wsDetailed.Range("B58").Value = 57.5792
$wsDetailed.Range("B59").Value = 68.16401999999999
$wsDetailed.Range("B60").Value = 68.66869

$wsDetailed.Range("B63").Value = 63.12887
$wsDetailed.Range("B64").Value = 30.4636
$wsDetailed.Range("B65").Value = 0.42209
$wsDetailed.Range("B66").Value = -2.60394
$wsDetailed.Range("B67").Value = -5.91381
$wsDetailed.Range("B68").Value = -6.77056
$wsDetailed.Range("B69").Value = -7.06268
$wsDetailed.Range("B70").Value = -9.5
$wsDetailed.Range("B71").Value = -9.84064

$wsDetailed.Range("B73").Value = -14
$wsDetailed.Range("B74").Value = -15.08844

$wsDetailed.Range("B76").Value = -21.85368
$wsDetailed.Range("B77").Value = -23.5

$wsDetailed.Range("B79").Value = -22.42548
$wsDetailed.Range("B80").Value = -23.06605
$wsDetailed.Range("B81").Value = -23.00352
$wsDetailed.Range("B82").Value = -7.34723
$wsDetailed.Range("B83").Value = -5.50985

$wsDetailed.Range("B85").Value = 47.19053
$wsDetailed.Range("B86").Value = 48.11341
$wsDetailed.Range("B87").Value = 64.02068
$wsDetailed.Range("B88").Value = 83.69806
$wsDetailed.Range("B89").Value = 105.40772
$wsDetailed.Range("B90").Value = 100.01
$wsDetailed.Range("B91").Value = 73.66
$wsDetailed.Range("B92").Value = 69.25239999999999

$wsDetailed.Range("B94").Value = 59.8861

$wsDetailed.Range("B96").Value = 64.8901
$wsDetailed.Range("B97").Value = 64.8901
